$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the data-row style (from row 2) to the new rows 14-22 before setting values
$ws.Range("A2:N2").Copy()
$ws.Range("A14:N22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2: Extra Trees
$ws.Cells.Item(2, 1).Value = "Extra Trees"
$ws.Cells.Item(2, 2).Value = 0.8761904761904762
$ws.Cells.Item(2, 3).Value = 0.9230769230769231
$ws.Cells.Item(2, 4).Value = 0.8484848484848485
$ws.Cells.Item(2, 5).Value = 0.8857808857808858
$ws.Cells.Item(2, 6).Value = 0.8811632811632812
$ws.Cells.Item(2, 7).Value = 0.782608695652174
$ws.Cells.Item(2, 8).Value = 0.9491525423728814
$ws.Cells.Item(2, 9).Value = 0.8658806190125277
$ws.Cells.Item(2, 10).Value = 0.8761904761904762
$ws.Cells.Item(2, 11).Value = 0.8470588235294118
$ws.Cells.Item(2, 12).Value = 0.896
$ws.Cells.Item(2, 13).Value = 0.8715294117647059
$ws.Cells.Item(2, 14).Value = 0.8745591036414566

# Row 3: Bagging
$ws.Cells.Item(3, 1).Value = "Bagging"
$ws.Cells.Item(3, 2).Value = 0.8380952380952381
$ws.Cells.Item(3, 3).Value = 0.8372093023255814
$ws.Cells.Item(3, 4).Value = 0.8387096774193549
$ws.Cells.Item(3, 5).Value = 0.8379594898724682
$ws.Cells.Item(3, 6).Value = 0.8380523702354161
$ws.Cells.Item(3, 7).Value = 0.782608695652174
$ws.Cells.Item(3, 8).Value = 0.8813559322033898
$ws.Cells.Item(3, 9).Value = 0.8319823139277819
$ws.Cells.Item(3, 10).Value = 0.8380952380952381
$ws.Cells.Item(3, 11).Value = 0.8089887640449438
$ws.Cells.Item(3, 12).Value = 0.859504132231405
$ws.Cells.Item(3, 13).Value = 0.8342464481381744
$ws.Cells.Item(3, 14).Value = 0.8373735899782887

# Row 4: KNN
$ws.Cells.Item(4, 1).Value = "KNN"
$ws.Cells.Item(4, 2).Value = 0.8095238095238095
$ws.Cells.Item(4, 3).Value = 0.8421052631578947
$ws.Cells.Item(4, 4).Value = 0.7910447761194029
$ws.Cells.Item(4, 5).Value = 0.8165750196386488
$ws.Cells.Item(4, 6).Value = 0.8134141323457897
$ws.Cells.Item(4, 7).Value = 0.6956521739130435
$ws.Cells.Item(4, 8).Value = 0.8983050847457628
$ws.Cells.Item(4, 9).Value = 0.7969786293294031
$ws.Cells.Item(4, 10).Value = 0.8095238095238095
$ws.Cells.Item(4, 11).Value = 0.761904761904762
$ws.Cells.Item(4, 12).Value = 0.8412698412698412
$ws.Cells.Item(4, 13).Value = 0.8015873015873016
$ws.Cells.Item(4, 14).Value = 0.8065003779289492

# Row 5: Random Forest
$ws.Cells.Item(5, 1).Value = "Random Forest"
$ws.Cells.Item(5, 2).Value = 0.8095238095238095
$ws.Cells.Item(5, 3).Value = 0.8611111111111112
$ws.Cells.Item(5, 4).Value = 0.782608695652174
$ws.Cells.Item(5, 5).Value = 0.8218599033816425
$ws.Cells.Item(5, 6).Value = 0.8170002300437084
$ws.Cells.Item(5, 7).Value = 0.6739130434782609
$ws.Cells.Item(5, 8).Value = 0.9152542372881356
$ws.Cells.Item(5, 9).Value = 0.7945836403831983
$ws.Cells.Item(5, 10).Value = 0.8095238095238095
$ws.Cells.Item(5, 11).Value = 0.7560975609756099
$ws.Cells.Item(5, 12).Value = 0.8437499999999999
$ws.Cells.Item(5, 13).Value = 0.7999237804878049
$ws.Cells.Item(5, 14).Value = 0.8053498838559814

# Row 6: MLP Classifier
$ws.Cells.Item(6, 1).Value = "MLP Classifier"
$ws.Cells.Item(6, 2).Value = 0.8
$ws.Cells.Item(6, 3).Value = 0.8205128205128205
$ws.Cells.Item(6, 4).Value = 0.7878787878787878
$ws.Cells.Item(6, 5).Value = 0.8041958041958042
$ws.Cells.Item(6, 6).Value = 0.8021756021756021
$ws.Cells.Item(6, 7).Value = 0.6956521739130435
$ws.Cells.Item(6, 8).Value = 0.8813559322033898
$ws.Cells.Item(6, 9).Value = 0.7885040530582166
$ws.Cells.Item(6, 10).Value = 0.8
$ws.Cells.Item(6, 11).Value = 0.7529411764705882
$ws.Cells.Item(6, 12).Value = 0.832
$ws.Cells.Item(6, 13).Value = 0.792470588235294
$ws.Cells.Item(6, 14).Value = 0.797364705882353

# Row 7: Gradient Boosting
$ws.Cells.Item(7, 1).Value = "Gradient Boosting"
$ws.Cells.Item(7, 2).Value = 0.7904761904761904
$ws.Cells.Item(7, 3).Value = 0.8157894736842105
$ws.Cells.Item(7, 4).Value = 0.7761194029850746
$ws.Cells.Item(7, 5).Value = 0.7959544383346426
$ws.Cells.Item(7, 6).Value = 0.7934986720532675
$ws.Cells.Item(7, 7).Value = 0.6739130434782609
$ws.Cells.Item(7, 8).Value = 0.8813559322033898
$ws.Cells.Item(7, 9).Value = 0.7776344878408253
$ws.Cells.Item(7, 10).Value = 0.7904761904761904
$ws.Cells.Item(7, 11).Value = 0.7380952380952381
$ws.Cells.Item(7, 12).Value = 0.8253968253968255
$ws.Cells.Item(7, 13).Value = 0.7817460317460319
$ws.Cells.Item(7, 14).Value = 0.7871504157218444

# Row 8: Decision Tree
$ws.Cells.Item(8, 1).Value = "Decision Tree"
$ws.Cells.Item(8, 2).Value = 0.780952380952381
$ws.Cells.Item(8, 3).Value = 0.7674418604651163
$ws.Cells.Item(8, 4).Value = 0.7903225806451613
$ws.Cells.Item(8, 5).Value = 0.7788822205551388
$ws.Cells.Item(8, 6).Value = 0.7802986460900939
$ws.Cells.Item(8, 7).Value = 0.717391304347826
$ws.Cells.Item(8, 8).Value = 0.8305084745762712
$ws.Cells.Item(8, 9).Value = 0.7739498894620487
$ws.Cells.Item(8, 10).Value = 0.780952380952381
$ws.Cells.Item(8, 11).Value = 0.7415730337078652
$ws.Cells.Item(8, 12).Value = 0.8099173553719008
$ws.Cells.Item(8, 13).Value = 0.775745194539883
$ws.Cells.Item(8, 14).Value = 0.7799760335000375

# Row 9 (AdaBoost) is unchanged from the original data
# Row 10: HistGradientBoosting
$ws.Cells.Item(10, 1).Value = "HistGradientBoosting"
$ws.Cells.Item(10, 2).Value = 0.7714285714285715
$ws.Cells.Item(10, 3).Value = 0.8055555555555556
$ws.Cells.Item(10, 4).Value = 0.7536231884057971
$ws.Cells.Item(10, 5).Value = 0.7795893719806763
$ws.Cells.Item(10, 6).Value = 0.7763745111571199
$ws.Cells.Item(10, 7).Value = 0.6304347826086957
$ws.Cells.Item(10, 8).Value = 0.8813559322033898
$ws.Cells.Item(10, 9).Value = 0.7558953574060427
$ws.Cells.Item(10, 10).Value = 0.7714285714285715
$ws.Cells.Item(10, 11).Value = 0.7073170731707318
$ws.Cells.Item(10, 12).Value = 0.8124999999999999
$ws.Cells.Item(10, 13).Value = 0.7599085365853658
$ws.Cells.Item(10, 14).Value = 0.7664198606271777

# Row 11: Quadratic Discriminant Analysis
$ws.Cells.Item(11, 1).Value = "Quadratic Discriminant Analysis"
$ws.Cells.Item(11, 2).Value = 0.7238095238095238
$ws.Cells.Item(11, 3).Value = 0.8148148148148148
$ws.Cells.Item(11, 4).Value = 0.6923076923076923
$ws.Cells.Item(11, 5).Value = 0.7535612535612535
$ws.Cells.Item(11, 6).Value = 0.7459774793108127
$ws.Cells.Item(11, 7).Value = 0.4782608695652174
$ws.Cells.Item(11, 8).Value = 0.9152542372881356
$ws.Cells.Item(11, 9).Value = 0.6967575534266766
$ws.Cells.Item(11, 10).Value = 0.7238095238095238
$ws.Cells.Item(11, 11).Value = 0.6027397260273973
$ws.Cells.Item(11, 12).Value = 0.7883211678832116
$ws.Cells.Item(11, 13).Value = 0.6955304469553045
$ws.Cells.Item(11, 14).Value = 0.7070188219273311

# Row 12: Linear SVC
$ws.Cells.Item(12, 1).Value = "Linear SVC"
$ws.Cells.Item(12, 2).Value = 0.7142857142857143
$ws.Cells.Item(12, 3).Value = 0.7666666666666667
$ws.Cells.Item(12, 4).Value = 0.6933333333333334
$ws.Cells.Item(12, 5).Value = 0.73
$ws.Cells.Item(12, 6).Value = 0.7254603174603175
$ws.Cells.Item(12, 7).Value = 0.5
$ws.Cells.Item(12, 8).Value = 0.8813559322033898
$ws.Cells.Item(12, 9).Value = 0.6906779661016949
$ws.Cells.Item(12, 10).Value = 0.7142857142857143
$ws.Cells.Item(12, 11).Value = 0.605263157894737
$ws.Cells.Item(12, 12).Value = 0.7761194029850748
$ws.Cells.Item(12, 13).Value = 0.6906912804399059
$ws.Cells.Item(12, 14).Value = 0.7012680956121649

# Row 13: Ridge Classifier
$ws.Cells.Item(13, 1).Value = "Ridge Classifier"
$ws.Cells.Item(13, 2).Value = 0.6952380952380952
$ws.Cells.Item(13, 3).Value = 0.7916666666666666
$ws.Cells.Item(13, 4).Value = 0.6666666666666666
$ws.Cells.Item(13, 5).Value = 0.7291666666666666
$ws.Cells.Item(13, 6).Value = 0.7214285714285714
$ws.Cells.Item(13, 7).Value = 0.4130434782608696
$ws.Cells.Item(13, 8).Value = 0.9152542372881356
$ws.Cells.Item(13, 9).Value = 0.6641488577745026
$ws.Cells.Item(13, 10).Value = 0.6952380952380952
$ws.Cells.Item(13, 11).Value = 0.5428571428571429
$ws.Cells.Item(13, 12).Value = 0.7714285714285714
$ws.Cells.Item(13, 13).Value = 0.6571428571428571
$ws.Cells.Item(13, 14).Value = 0.6712925170068028

# Row 14: SGD Classifier
$ws.Cells.Item(14, 1).Value = "SGD Classifier"
$ws.Cells.Item(14, 2).Value = 0.6952380952380952
$ws.Cells.Item(14, 3).Value = 0.6590909090909091
$ws.Cells.Item(14, 4).Value = 0.7213114754098361
$ws.Cells.Item(14, 5).Value = 0.6902011922503726
$ws.Cells.Item(14, 6).Value = 0.6940529415939252
$ws.Cells.Item(14, 7).Value = 0.6304347826086957
$ws.Cells.Item(14, 8).Value = 0.7457627118644068
$ws.Cells.Item(14, 9).Value = 0.6880987472365512
$ws.Cells.Item(14, 10).Value = 0.6952380952380952
$ws.Cells.Item(14, 11).Value = 0.6444444444444444
$ws.Cells.Item(14, 12).Value = 0.7333333333333334
$ws.Cells.Item(14, 13).Value = 0.6888888888888889
$ws.Cells.Item(14, 14).Value = 0.6943915343915344

# Row 15: Passive Aggressive
$ws.Cells.Item(15, 1).Value = "Passive Aggressive"
$ws.Cells.Item(15, 2).Value = 0.6952380952380952
$ws.Cells.Item(15, 3).Value = 0.7692307692307693
$ws.Cells.Item(15, 4).Value = 0.6708860759493671
$ws.Cells.Item(15, 5).Value = 0.7200584225900681
$ws.Cells.Item(15, 6).Value = 0.7139704177678861
$ws.Cells.Item(15, 7).Value = 0.4347826086956522
$ws.Cells.Item(15, 8).Value = 0.8983050847457628
$ws.Cells.Item(15, 9).Value = 0.6665438467207074
$ws.Cells.Item(15, 10).Value = 0.6952380952380952
$ws.Cells.Item(15, 11).Value = 0.5555555555555555
$ws.Cells.Item(15, 12).Value = 0.7681159420289855
$ws.Cells.Item(15, 13).Value = 0.6618357487922705
$ws.Cells.Item(15, 14).Value = 0.6749942489072923

# Row 16: SVM
$ws.Cells.Item(16, 1).Value = "SVM"
$ws.Cells.Item(16, 2).Value = 0.6761904761904762
$ws.Cells.Item(16, 3).Value = 0.7727272727272727
$ws.Cells.Item(16, 4).Value = 0.6506024096385542
$ws.Cells.Item(16, 5).Value = 0.7116648411829134
$ws.Cells.Item(16, 6).Value = 0.7041047306107546
$ws.Cells.Item(16, 7).Value = 0.3695652173913043
$ws.Cells.Item(16, 8).Value = 0.9152542372881356
$ws.Cells.Item(16, 9).Value = 0.64240972733972
$ws.Cells.Item(16, 10).Value = 0.6761904761904762
$ws.Cells.Item(16, 11).Value = 0.4999999999999999
$ws.Cells.Item(16, 12).Value = 0.7605633802816901
$ws.Cells.Item(16, 13).Value = 0.630281690140845
$ws.Cells.Item(16, 14).Value = 0.646411804158283

# Row 17: Logistic Regression
$ws.Cells.Item(17, 1).Value = "Logistic Regression"
$ws.Cells.Item(17, 2).Value = 0.6666666666666666
$ws.Cells.Item(17, 3).Value = 0.7894736842105263
$ws.Cells.Item(17, 4).Value = 0.6395348837209303
$ws.Cells.Item(17, 5).Value = 0.7145042839657283
$ws.Cells.Item(17, 6).Value = 0.7052223582211343
$ws.Cells.Item(17, 7).Value = 0.3260869565217391
$ws.Cells.Item(17, 8).Value = 0.9322033898305084
$ws.Cells.Item(17, 9).Value = 0.6291451731761237
$ws.Cells.Item(17, 10).Value = 0.6666666666666666
$ws.Cells.Item(17, 11).Value = 0.4615384615384616
$ws.Cells.Item(17, 12).Value = 0.7586206896551724
$ws.Cells.Item(17, 13).Value = 0.6100795755968169
$ws.Cells.Item(17, 14).Value = 0.6284703801945181

# Row 18: Linear Discriminant Analysis
$ws.Cells.Item(18, 1).Value = "Linear Discriminant Analysis"
$ws.Cells.Item(18, 2).Value = 0.6476190476190476
$ws.Cells.Item(18, 3).Value = 0.6666666666666666
$ws.Cells.Item(18, 4).Value = 0.6410256410256411
$ws.Cells.Item(18, 5).Value = 0.6538461538461539
$ws.Cells.Item(18, 6).Value = 0.6522588522588523
$ws.Cells.Item(18, 7).Value = 0.391304347826087
$ws.Cells.Item(18, 8).Value = 0.847457627118644
$ws.Cells.Item(18, 9).Value = 0.6193809874723655
$ws.Cells.Item(18, 10).Value = 0.6476190476190476
$ws.Cells.Item(18, 11).Value = 0.4931506849315068
$ws.Cells.Item(18, 12).Value = 0.7299270072992702
$ws.Cells.Item(18, 13).Value = 0.6115388461153886
$ws.Cells.Item(18, 14).Value = 0.6261964279762501

# Row 19: Ridge Classifier CV
$ws.Cells.Item(19, 1).Value = "Ridge Classifier CV"
$ws.Cells.Item(19, 2).Value = 0.6476190476190476
$ws.Cells.Item(19, 3).Value = 0.7647058823529411
$ws.Cells.Item(19, 4).Value = 0.625
$ws.Cells.Item(19, 5).Value = 0.6948529411764706
$ws.Cells.Item(19, 6).Value = 0.6862044817927171
$ws.Cells.Item(19, 7).Value = 0.2826086956521739
$ws.Cells.Item(19, 8).Value = 0.9322033898305084
$ws.Cells.Item(19, 9).Value = 0.6074060427413411
$ws.Cells.Item(19, 10).Value = 0.6476190476190476
$ws.Cells.Item(19, 11).Value = 0.4126984126984127
$ws.Cells.Item(19, 12).Value = 0.7482993197278912
$ws.Cells.Item(19, 13).Value = 0.5804988662131519
$ws.Cells.Item(19, 14).Value = 0.6012741604578339

# Row 20: Nearest Centroid
$ws.Cells.Item(20, 1).Value = "Nearest Centroid"
$ws.Cells.Item(20, 2).Value = 0.6285714285714286
$ws.Cells.Item(20, 3).Value = 0.7058823529411765
$ws.Cells.Item(20, 4).Value = 0.6136363636363636
$ws.Cells.Item(20, 5).Value = 0.6597593582887701
$ws.Cells.Item(20, 6).Value = 0.6540488922841865
$ws.Cells.Item(20, 7).Value = 0.2608695652173913
$ws.Cells.Item(20, 8).Value = 0.9152542372881356
$ws.Cells.Item(20, 9).Value = 0.5880619012527635
$ws.Cells.Item(20, 10).Value = 0.6285714285714286
$ws.Cells.Item(20, 11).Value = 0.380952380952381
$ws.Cells.Item(20, 12).Value = 0.7346938775510204
$ws.Cells.Item(20, 13).Value = 0.5578231292517007
$ws.Cells.Item(20, 14).Value = 0.5797214123744736

# Row 21: Perceptron
$ws.Cells.Item(21, 1).Value = "Perceptron"
$ws.Cells.Item(21, 2).Value = 0.6190476190476191
$ws.Cells.Item(21, 3).Value = 0.5535714285714286
$ws.Cells.Item(21, 4).Value = 0.6938775510204082
$ws.Cells.Item(21, 5).Value = 0.6237244897959184
$ws.Cells.Item(21, 6).Value = 0.6324101068999027
$ws.Cells.Item(21, 7).Value = 0.6739130434782609
$ws.Cells.Item(21, 8).Value = 0.576271186440678
$ws.Cells.Item(21, 9).Value = 0.6250921149594695
$ws.Cells.Item(21, 10).Value = 0.6190476190476191
$ws.Cells.Item(21, 11).Value = 0.6078431372549019
$ws.Cells.Item(21, 12).Value = 0.6296296296296295
$ws.Cells.Item(21, 13).Value = 0.6187363834422657
$ws.Cells.Item(21, 14).Value = 0.6200850710654632

# Row 22: Naive Bayes
$ws.Cells.Item(22, 1).Value = "Naive Bayes"
$ws.Cells.Item(22, 2).Value = 0.5904761904761905
$ws.Cells.Item(22, 3).Value = 0.7142857142857143
$ws.Cells.Item(22, 4).Value = 0.5816326530612245
$ws.Cells.Item(22, 5).Value = 0.6479591836734694
$ws.Cells.Item(22, 6).Value = 0.6397473275024296
$ws.Cells.Item(22, 7).Value = 0.108695652173913
$ws.Cells.Item(22, 8).Value = 0.9661016949152542
$ws.Cells.Item(22, 9).Value = 0.5373986735445836
$ws.Cells.Item(22, 10).Value = 0.5904761904761905
$ws.Cells.Item(22, 11).Value = 0.1886792452830189
$ws.Cells.Item(22, 12).Value = 0.7261146496815286
$ws.Cells.Item(22, 13).Value = 0.4573969474822737
$ws.Cells.Item(22, 14).Value = 0.4906667582307529
